# Weekly price update: insert two new price records (new week of data,
# 2021-12-27) at the top of the "Caigua" detail table, pushing the
# existing rows (36..86) down by two positions. The two oldest rows that
# fall past the end of the original range end up appended as new rows
# 87 and 88 (handled automatically by the row insert / shift-down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 36, shifting
# everything below down by two rows (36->38 ... 86->88).
$ws.Rows("36:37").Insert()

# New row 36: "Primera" quality, week of 2021-12-27
$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = "12/27/2021"
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = 100112036
$ws.Range("G36").Value = "Caigua"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 160
$ws.Range("K36").Value = 5000
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = 5500
$ws.Range("N36").Value = '$/caja 20 kilos'
$ws.Range("O36").Value = "Región de Arica y Parinacota"
$ws.Range("P36").Value = 275
$ws.Range("Q36").Value = 20
$ws.Range("R36").Value = "Hortaliza"

# New row 37: "Segunda" quality, week of 2021-12-27
$ws.Range("A37").Value = 1
$ws.Range("B37").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C37").Value = "Arica y Parinacota"
$ws.Range("D37").Value = "12/27/2021"
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = 100112036
$ws.Range("G37").Value = "Caigua"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Segunda"
$ws.Range("J37").Value = 120
$ws.Range("K37").Value = 4500
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = 4750
$ws.Range("N37").Value = '$/caja 20 kilos'
$ws.Range("O37").Value = "Región de Arica y Parinacota"
$ws.Range("P37").Value = 238
$ws.Range("Q37").Value = 20
$ws.Range("R37").Value = "Hortaliza"
